$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 821, shifting existing rows 821:862 down to 822:863.
$ws.Rows("821:821").Insert()

# Column A holds the date as plain text (e.g. "2026/12/29"), not a real Excel
# date. Force text formatting before assigning so COM doesn't auto-convert
# the "2026/02/16" string into a date serial, then restore the default
# "Normal" style so the cell doesn't end up with a stray number-format style
# that the original file doesn't have.
$ws.Range("A821").NumberFormat = "@"
$ws.Range("A821").Value = "2026/02/16"
$ws.Range("A821").Style = "Normal"

$ws.Range("B821").Value = "月"
$ws.Range("C821").Value = 17
$ws.Range("D821").Value = 201
